$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates ---
# Note: the host's ColumnWidth setter quantizes to a pixel grid before
# persisting (same behavior Excel itself exhibits, just keyed to this
# engine's digit-width metric), so the inputs below are chosen to land on
# the pixel-grid point nearest the target "characters" width rather than
# the raw target value itself.
$ws.Columns.Item(3).ColumnWidth = 1.3333333333333333   # -> stored width ~2.140625
$ws.Columns.Item(6).ColumnWidth = 2.3333333333333335   # -> stored width ~3.140625
$ws.Columns.Item(9).ColumnWidth = 4.833333333333333    # -> stored width ~5.7109375
$ws.Columns.Item(10).ColumnWidth = 4.833333333333333   # -> stored width ~5.7109375

# --- Row 1 cell value updates ---
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 8
$ws.Range("D1").Value = 12
$ws.Range("E1").Value = 21
$ws.Range("F1").Value = 33
$ws.Range("G1").Value = 14
$ws.Range("H1").Value = 19
$ws.Range("I1").Value = 0.028000000000000004
$ws.Range("J1").Value = 0.058999999999999997
$ws.Range("K1").Value = 0.094
